$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.139845967292786
$ws.Range("B1").Value = 2.415846109390259
$ws.Range("C1").Value = 5.098052024841309
$ws.Range("D1").Value = 2.226372003555298
$ws.Range("E1").Value = 1.253384351730347
